# Thailand Premier League - base update (11-04-2024 23:56)
# The source feed re-sorted a handful of same-date/same-round fixture pairs,
# which swapped the row order for those matches. Swap back the full data
# (every column except the row's running "id" in column A) between each
# affected pair of rows so the match details line up with the right id again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, [int]$row1, [int]$row2, [int]$firstCol, [int]$lastCol)

    $row1Values = @()
    $row2Values = @()

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $row1Values += , $ws.Cells.Item($row1, $col).Value()
        $row2Values += , $ws.Cells.Item($row2, $col).Value()
    }

    $i = 0
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row1, $col).Value = $row2Values[$i]
        $ws.Cells.Item($row2, $col).Value = $row1Values[$i]
        $i++
    }
}

# Columns B (2) through AC (29) hold the match data; column A (id) is left untouched.
Swap-RowData $ws 85 86 2 29
Swap-RowData $ws 117 118 2 29
Swap-RowData $ws 179 180 2 29
